$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ----------------------------------------------------------------------------
# 1) Shape "Textfeld 45" (id 46): fix typo "Battry" -> "Battery" and widen /
#    re-position the box slightly (off.x 3931219 -> 3893581, ext.cx 999954 ->
#    1075231; off.y and ext.cy stay the same).
# ----------------------------------------------------------------------------
$battery = $s.Shapes.Item(21)
# Target only the "Battry" run (first 6 characters) so the fix-up stays a
# single run, just like the rest of the paragraph ("DC/DC ", "switch", ...).
$battery.TextFrame.TextRange.Characters(1, 6).Text = "Battery"
$battery.Left  = 306.58122110236224
$battery.Width = 84.66389826771653

# ----------------------------------------------------------------------------
# 2) Shape "Textfeld 46" (id 47): nudge down a little (off.y 469628 -> 516521).
# ----------------------------------------------------------------------------
$dcdc = $s.Shapes.Item(22)
$dcdc.Top = 40.670984881889765

# ----------------------------------------------------------------------------
# 3) New shape "Textfeld 47" (id 48) labelling the Speed pin, placed next to
#    the other rotated pin labels.
#
# The host runtime hands out shape ids from a monotonically increasing
# counter that is seeded by every Shapes.Add*/Duplicate call ever made in
# the session (ids already present in the file are skipped, but freeing an
# id by deleting a shape does not let it be re-used). The real authored file
# has this textbox at id 48, i.e. the 26th shape ever created in the file's
# history (after ids 2,3,5,23-34,36-45 get burned by the 25 preceding
# inserts). Burn through 25 throw-away ids first so the textbox we actually
# keep lands on id 48, matching the source file exactly.
# ----------------------------------------------------------------------------
for ($i = 1; $i -le 25; $i++) {
    $scratch = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $scratch.Delete()
}

# Duplicate the sibling pin-label textbox so the new shape inherits identical
# formatting (noFill, wrap="none"/spAutoFit body, de-DE/12pt run properties).
$speedRange = $dcdc.Duplicate()
$speed = $speedRange.Item(1)
$speed.Name = "Textfeld 47"
$speed.TextFrame.TextRange.Text = "Speed"
$speed.Left   = 359.8285439370079
$speed.Top    = 451.60909511811025
$speed.Width  = 44.83366204724409
$speed.Height = 21.810984881889762
$speed.Rotation = 270
